# Prosaanforderungen.xlsx -- "erledigt; kosmetische Aenderungen"
# Fills in Tabelle1 (sheet1) with the full list of 20 prose requirements,
# re-centers the "Nr." index column, and relocates the illustration from
# Tabelle1 to Tabelle2.

$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item(1)
$ws2 = $wb.Worksheets.Item(2)

# ---------------------------------------------------------------------
# 1) Header row -- write B1 ("Text") before A1 ("Nr.") so the shared
#    string table is rebuilt in that exact order (index 0 / 1).
# ---------------------------------------------------------------------
$ws1.Range("B1").Value = "Text"
$ws1.Range("A1").Value = "Nr."
$ws1.Range("A1:B1").HorizontalAlignment = -4108   # xlCenter (keeps existing bold/center style)

# ---------------------------------------------------------------------
# 2) Body rows -- insertion order drives shared-string index assignment,
#    so the *first* eight unique sentences below become indices 2..9 and
#    the remaining twelve become indices 10..21, independent of which
#    sheet row they eventually live in.
# ---------------------------------------------------------------------
$bodyByRow = @{
  14 = "Die persistenten Daten müssen auf einer`nPostgreSQL-DB via JDBC abgespeichert werden."
  15 = "Mit dem Tool SonarQube müssen`ntechnische Schulden und Bad Code Smells`nAnalysiert und visualisiert werden."
  16 = "Die Datenbank muss mit Roundtrip-Tests getestet werden."
  17 = "Die Software muss via SEPP deployed werden."
  18 = "Der externe Zugriff zum deployen der Software via SEPP`nmuss via OpenVPN erfolgen."
  19 = "Benachrichtigungen über Software-Fehler sollen per E-Mail gesendet werden."
  20 = "In der Produktivumgebung soll die E-Mail-Benachrichtigung`nim Fehlerfall deaktiviert werden können."
  21 = "Die Datenbankimplementierung soll mit Mockito getestet werden."
  2  = "Der studentische Benutzer soll sich mit seiner Linux-Kennung einloggen können."
  3  = "Der studentische Benutzer soll seinen Account löschen können, um keine weiteren Benachrichtigungen zu bekommen."
  4  = "Unternehmerische Benutzer sollen sich registrieren können."
  5  = "Der studentische Benutzer soll sich ein Profil mit persönlichen Daten und fachlichen Interessen erstellen können."
  7  = "Der studentische Benutzer soll sich nach seinen Qualifikationen, Stärken, befristeten Stellen und Beschäftigungsart filtern können, um passende Stellen zu finden."
  10 = "Der studentische Benutzer soll anhand seiner fachlichen Interessen zu passenden Stellen informiert werden."
  11 = "Der studentische Benutzer soll Unternehmen abonnieren können, um über Neuigkeiten & Stellen informiert zu werden."
  12 = "Der Benutzer benötigt eine passwortgeschützte Anmeldung."
  13 = "Der Benutzer soll eine auf Smartphones angepasste Ansicht des Systems haben."
  6  = "Unternehmerische Benutzer sollen Informationen über ihr Unternehmen bereitstellen können."
  8  = "Unternehmerische Benutzer sollen freie Stellen bereitstellen können."
  9  = "Unternehmerische Benutzer sollen einen Ansprechpartner hinterlegen können."
}

$insertOrder = @(14,15,16,17,18,19,20,21,2,3,4,5,7,10,11,12,13,6,8,9)
foreach ($row in $insertOrder) {
  $ws1.Cells.Item($row, 2).Value = $bodyByRow[$row]
}

# ---------------------------------------------------------------------
# 3) Column A -- sequential numbers 1..20, center aligned (new style).
# ---------------------------------------------------------------------
for ($row = 2; $row -le 21; $row++) {
  $ws1.Cells.Item($row, 1).Value = $row - 1
}
$ws1.Range("A2:A21").HorizontalAlignment = -4108   # xlCenter -> creates the new default-font centered style

# ---------------------------------------------------------------------
# 4) Column / selection cosmetics.
# ---------------------------------------------------------------------
$ws1.Columns.Item(1).ColumnWidth = 3.6640625
$ws1.Range("B1").Select()

# ---------------------------------------------------------------------
# 5) Move the illustration from Tabelle1 to Tabelle2, anchored at A1,
#    same physical size as before.
# ---------------------------------------------------------------------
$ws1.Shapes.Item(1).Delete()
$pic = $ws2.Shapes.AddPicture("/tmp/work/image1.png", 0, 1, 0, 0, 478.716062992126, 328.2)
